# Auto-generated edit script applying numeric corrections to Leve profit-tracking sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value2 = 981.0625
$ws.Range("I15").Value2 = 981.0625
$ws.Range("K15").Value2 = 2943.1875
$ws.Range("M15").Value2 = -2774.1875

# Row 107
$ws.Range("H107").Value2 = 3264
$ws.Range("I107").Value2 = 3264
$ws.Range("K107").Value2 = 3264
$ws.Range("M107").Value2 = -1344

# Row 113
$ws.Range("H113").Value2 = 4271.273
$ws.Range("I113").Value2 = 4164.8887
$ws.Range("J113").Value2 = 4750
$ws.Range("K113").Value2 = 4164.8887
$ws.Range("L113").Value2 = 4750
$ws.Range("M113").Value2 = -910.8887000000004
$ws.Range("N113").Value2 = -11258

# Row 116
$ws.Range("H116").Value2 = 35183
$ws.Range("J116").Value2 = 63724.75
$ws.Range("L116").Value2 = 63724.75
$ws.Range("N116").Value2 = -70608.75

# Row 132
$ws.Range("H132").Value2 = 20344.637
$ws.Range("I132").Value2 = 24426.777
$ws.Range("K132").Value2 = 73280.33099999999
$ws.Range("M132").Value2 = -70750.33099999999

# Row 136
$ws.Range("H136").Value2 = 90000
$ws.Range("J136").Value2 = 90000
$ws.Range("L136").Value2 = 90000
$ws.Range("N136").Value2 = -100200

# Row 137
$ws.Range("H137").Value2 = 7696.45
$ws.Range("I137").Value2 = 2838.3684
$ws.Range("K137").Value2 = 8515.1052
$ws.Range("M137").Value2 = -5965.1052

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value2 = 2345.5557
$ws.Range("I45").Value2 = 2010.1666
$ws.Range("J45").Value2 = 3016.3333
$ws.Range("K45").Value2 = 2010.1666
$ws.Range("L45").Value2 = 3016.3333
$ws.Range("M45").Value2 = -1633.1666
$ws.Range("N45").Value2 = -3770.3333

# Row 61
$ws.Range("H61").Value2 = 2067.842
$ws.Range("I61").Value2 = 2005.2354
$ws.Range("K61").Value2 = 2005.2354
$ws.Range("M61").Value2 = -1793.2354

# Row 122
$ws.Range("H122").Value2 = 2396.5
$ws.Range("I122").Value2 = 1904.5834
$ws.Range("J122").Value2 = 8299.5
$ws.Range("K122").Value2 = 5713.7502
$ws.Range("L122").Value2 = 24898.5
$ws.Range("M122").Value2 = -3263.7502
$ws.Range("N122").Value2 = -29798.5

# Row 132
$ws.Range("H132").Value2 = 7956.731
$ws.Range("I132").Value2 = 5898.909
$ws.Range("J132").Value2 = 9465.799999999999
$ws.Range("K132").Value2 = 17696.727
$ws.Range("L132").Value2 = 28397.4
$ws.Range("M132").Value2 = -15166.727
$ws.Range("N132").Value2 = -33457.39999999999

# Row 136
$ws.Range("H136").Value2 = 2067.842
$ws.Range("I136").Value2 = 2005.2354
$ws.Range("K136").Value2 = 6015.706200000001
$ws.Range("M136").Value2 = -3465.706200000001

# Row 137
$ws.Range("H137").Value2 = 90000
$ws.Range("J137").Value2 = 90000
$ws.Range("L137").Value2 = 90000
$ws.Range("N137").Value2 = -100200

# Row 141
$ws.Range("H141").Value2 = 0
$ws.Range("J141").Value2 = 0
$ws.Range("N141").Value2 = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value2 = 75771.42999999999
$ws.Range("I20").Value2 = 116036.22
$ws.Range("K20").Value2 = 116036.22
$ws.Range("M20").Value2 = -115789.22

# Row 26
$ws.Range("H26").Value2 = 24948.334
$ws.Range("I26").Value2 = 16938
$ws.Range("J26").Value2 = 65000
$ws.Range("K26").Value2 = 16938
$ws.Range("L26").Value2 = 65000
$ws.Range("M26").Value2 = -16646
$ws.Range("N26").Value2 = -65584

# Row 94
$ws.Range("H94").Value2 = 8354.9
$ws.Range("I94").Value2 = 7616.5557
$ws.Range("J94").Value2 = 15000
$ws.Range("K94").Value2 = 7616.5557
$ws.Range("L94").Value2 = 15000
$ws.Range("M94").Value2 = -7165.5557
$ws.Range("N94").Value2 = -15902

# Row 99
$ws.Range("H99").Value2 = 13876.25
$ws.Range("I99").Value2 = 13876.25
$ws.Range("K99").Value2 = 13876.25
$ws.Range("M99").Value2 = -12378.25

# Row 134
$ws.Range("H134").Value2 = 1574.3846
$ws.Range("I134").Value2 = 946.8
$ws.Range("J134").Value2 = 3666.3333
$ws.Range("K134").Value2 = 2840.4
$ws.Range("L134").Value2 = 10998.9999
$ws.Range("M134").Value2 = -305.3999999999996
$ws.Range("N134").Value2 = -16068.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 3416
$ws.Range("J31").Value2 = 4600
$ws.Range("L31").Value2 = 4600
$ws.Range("N31").Value2 = -5190

# Row 34
$ws.Range("H34").Value2 = 3416
$ws.Range("J34").Value2 = 4600
$ws.Range("L34").Value2 = 4600
$ws.Range("N34").Value2 = -5004

# Row 58
$ws.Range("H58").Value2 = 5751.8125
$ws.Range("J58").Value2 = 6265
$ws.Range("L58").Value2 = 6265
$ws.Range("N58").Value2 = -6671

# Row 105
$ws.Range("H105").Value2 = 26141.25
$ws.Range("I105").Value2 = 26141.25
$ws.Range("J105").Value2 = 0
$ws.Range("K105").Value2 = 26141.25
$ws.Range("L105").Value2 = 0
$ws.Range("N105").Value2 = -24394.25
$ws.Range("N105").ClearContents()

# Row 107
$ws.Range("H107").Value2 = 1648.7273
$ws.Range("I107").Value2 = 1506.0834
$ws.Range("J107").Value2 = 1819.9
$ws.Range("K107").Value2 = 1506.0834
$ws.Range("L107").Value2 = 1819.9
$ws.Range("M107").Value2 = 413.9166
$ws.Range("N107").Value2 = -5659.9

# Row 132
$ws.Range("H132").Value2 = 3182.8696
$ws.Range("I132").Value2 = 3040.3
$ws.Range("K132").Value2 = 9120.900000000001
$ws.Range("M132").Value2 = -6590.900000000001

# Row 134
$ws.Range("H134").Value2 = 2763.75
$ws.Range("I134").Value2 = 2167.2856
$ws.Range("K134").Value2 = 6501.8568
$ws.Range("M134").Value2 = -3966.8568

# Row 136
$ws.Range("H136").Value2 = 5751.8125
$ws.Range("J136").Value2 = 6265
$ws.Range("L136").Value2 = 18795
$ws.Range("N136").Value2 = -23895

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value2 = 3324.25
$ws.Range("I68").Value2 = 1198.4286
$ws.Range("J68").Value2 = 4977.6665
$ws.Range("K68").Value2 = 3595.2858
$ws.Range("L68").Value2 = 14932.9995
$ws.Range("M68").Value2 = -2784.2858
$ws.Range("N68").Value2 = -16554.9995

# Row 71
$ws.Range("H71").Value2 = 3324.25
$ws.Range("I71").Value2 = 1198.4286
$ws.Range("J71").Value2 = 4977.6665
$ws.Range("K71").Value2 = 10785.8574
$ws.Range("L71").Value2 = 44798.9985
$ws.Range("M71").Value2 = -6729.857399999999
$ws.Range("N71").Value2 = -52910.9985

# Row 104
$ws.Range("H104").Value2 = 5116
$ws.Range("I104").Value2 = 6000
$ws.Range("K104").Value2 = 18000
$ws.Range("M104").Value2 = -15379

# Row 122
$ws.Range("H122").Value2 = 1153270.6
$ws.Range("J122").Value2 = 1358.9166
$ws.Range("L122").Value2 = 12230.2494
$ws.Range("N122").Value2 = -17130.2494

# Row 134
$ws.Range("H134").Value2 = 7432.9414
$ws.Range("I134").Value2 = 3978.182
$ws.Range("K134").Value2 = 11934.546
$ws.Range("M134").Value2 = -6864.545999999998

# Row 136
$ws.Range("H136").Value2 = 10547.462
$ws.Range("J136").Value2 = 13614.875
$ws.Range("L136").Value2 = 40844.625
$ws.Range("N136").Value2 = -51044.625

# Row 141
$ws.Range("H141").Value2 = 9067.3125
$ws.Range("I141").Value2 = 4507.7
$ws.Range("J141").Value2 = 16666.666
$ws.Range("K141").Value2 = 13523.1
$ws.Range("L141").Value2 = 49999.99800000001
$ws.Range("M141").Value2 = -8343.099999999999
$ws.Range("N141").Value2 = -60359.99800000001

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value2 = 23082.666
$ws.Range("J57").Value2 = 23082.666
$ws.Range("L57").Value2 = 23082.666
$ws.Range("N57").Value2 = -24722.666

# Row 122
$ws.Range("H122").Value2 = 1908.6364
$ws.Range("I122").Value2 = 1519.6666
$ws.Range("K122").Value2 = 4558.9998
$ws.Range("M122").Value2 = -2108.9998

# Row 140
$ws.Range("H140").Value2 = 97390
$ws.Range("J140").Value2 = 97390
$ws.Range("L140").Value2 = 97390
$ws.Range("N140").Value2 = -107750

# Row 141
$ws.Range("H141").Value2 = 77733.125
$ws.Range("J141").Value2 = 77733.125
$ws.Range("L141").Value2 = 77733.125
$ws.Range("N141").Value2 = -88093.125

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value2 = 2789.3809
$ws.Range("I22").Value2 = 2437.923
$ws.Range("J22").Value2 = 3360.5
$ws.Range("K22").Value2 = 2437.923
$ws.Range("L22").Value2 = 3360.5
$ws.Range("M22").Value2 = -2142.923
$ws.Range("N22").Value2 = -3950.5

# Row 27
$ws.Range("H27").Value2 = 2789.3809
$ws.Range("I27").Value2 = 2437.923
$ws.Range("J27").Value2 = 3360.5
$ws.Range("K27").Value2 = 2437.923
$ws.Range("L27").Value2 = 3360.5
$ws.Range("M27").Value2 = -2330.923
$ws.Range("N27").Value2 = -3574.5

# Row 46
$ws.Range("H46").Value2 = 4011.1304
$ws.Range("I46").Value2 = 1324.1428
$ws.Range("J46").Value2 = 5186.6875
$ws.Range("K46").Value2 = 1324.1428
$ws.Range("L46").Value2 = 5186.6875
$ws.Range("M46").Value2 = -1136.1428
$ws.Range("N46").Value2 = -5562.6875

# Row 100
$ws.Range("H100").Value2 = 2723.7144
$ws.Range("I100").Value2 = 2343.6667
$ws.Range("J100").Value2 = 5004
$ws.Range("K100").Value2 = 2343.6667
$ws.Range("L100").Value2 = 5004
$ws.Range("M100").Value2 = -1802.6667
$ws.Range("N100").Value2 = -6086

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value2 = 2339.0667
$ws.Range("I126").Value2 = 1798
$ws.Range("J126").Value2 = 2957.4285
$ws.Range("K126").Value2 = 5394
$ws.Range("L126").Value2 = 8872.2855
$ws.Range("M126").Value2 = -2924
$ws.Range("N126").Value2 = -13812.2855

# Row 132
$ws.Range("H132").Value2 = 2033.6364
$ws.Range("I132").Value2 = 1590.8857
$ws.Range("J132").Value2 = 3755.4443
$ws.Range("K132").Value2 = 4772.6571
$ws.Range("L132").Value2 = 11266.3329
$ws.Range("M132").Value2 = -2242.6571
$ws.Range("N132").Value2 = -16326.3329

# Row 136
$ws.Range("H136").Value2 = 1426.2413
$ws.Range("I136").Value2 = 1362.6666
$ws.Range("J136").Value2 = 2284.5
$ws.Range("K136").Value2 = 4087.9998
$ws.Range("L136").Value2 = 6853.5
$ws.Range("M136").Value2 = -1537.9998
$ws.Range("N136").Value2 = -11953.5
